$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $rng = $d.Range($full.Start, $full.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

$p5 = @'
<w:p><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Voorlopige basisuitspraak:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">de vooral onbekende keuzes die door de ECB onder Lagarde zijn gemaakt in de periode 2019-2020, zijn te laat en onverantwoord gemaakt </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">De huidige inflatie is hiervan het resultaat, met als gevolg dat de ECB zijn renteverhogingen op een snel tempo heeft moeten invoeren </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Prijzen in de supermarkt zijn voor burgers gigantisch gestegen met zelfs unieke fenomenen als </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>krimpflatie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> als gevolg.</w:t></w:r></w:p>
'@
Set-ParagraphXml 5 $p5

$p8 = @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Subthema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> 1</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Het naast elkaar houden van twee supermarkt mandjes eentje uit 2023 en de ander uit 2019. Dit verschil van 4 jaar heeft grote veranderingen met zich mee gebracht voor het winkelmandje voor de weekboodschappen van de gemiddelde Nederlander</w:t></w:r></w:p>
'@
Set-ParagraphXml 8 $p8

$p10 = @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Subthema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> 2</w:t></w:r></w:p>
'@
Set-ParagraphXml 10 $p10

$p13 = @'
<w:p><w:r><w:t xml:space="preserve">Dat wil ik bespreken met econoom Edin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mujagic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wie voor verschillende instanties het </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Macroeconomische</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nieuws brengt.</w:t></w:r></w:p>
'@
Set-ParagraphXml 13 $p13

$p16 = @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Subthema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> 3</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>Wat heeft de ECB t</w:t></w:r><w:r><w:t xml:space="preserve">oen gedaan? De ECB heeft hierop gereageerd door in enorme aantallen Eurobiljetten te printen. En hierbij de Europese economie te steunen die bezig waren met massale stimulus </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paketten</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> om stilstaande </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>economieen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> te stutten.</w:t></w:r></w:p>
'@
Set-ParagraphXml 16 $p16

$p18 = @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Subthema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>4</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Resultaat van beleid ECB en krimp- &amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>graaiflatie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>deze punten hebben geleid tot de huidige situatie voor het afgelopen jaar. De prijzen van consumentengoederen daalt maar niet en er treden nu ook nieuwe economische fenomenen op zoals ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>graaiflatie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’ en ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>krimpflatie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’. Hoe kunnen we dit herkennen en hoe gaat dit in de toekomst eruit zien?</w:t></w:r></w:p>
'@
Set-ParagraphXml 18 $p18

Write-Output "done"